$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the mmrsurv (governance delta) sensitivity value from 0.5 to 0.1
# and populate the resulting recomputed indicator values for each affected row.

# Row 4 (afg / 2006)
$ws.Range("E4").Value = 0.1
$ws.Range("F4").Value = 42.060589412729861
$ws.Range("J4").Value = 59.805837819685323
$ws.Range("K4").Value = 47.343431044127705

# Row 7 (alb / 2011)
$ws.Range("E7").Value = 0.1
$ws.Range("G7").Value = 66.263254198720958
$ws.Range("N7").Value = 66.265000000000001

# Row 10 (alb / 2016)
$ws.Range("E10").Value = 0.1
$ws.Range("H10").Value = 99.206287487457246
$ws.Range("I10").Value = 42.668358791697457
$ws.Range("N10").Value = 42.670999999999999

# Row 13 (ago / 2017)
$ws.Range("E13").Value = 0.1
$ws.Range("L13").Value = 93.772924571552139
$ws.Range("M13").Value = 99.903619964787666

# Move the active selection to H10, matching the saved view state
$ws.Range("H10").Select()
